$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, pushing existing rows 16-93 down to 17-94.
$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with the new weekly data record.
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 45114
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112013
$ws.Range("G16").Value = "Alcachofa"
$ws.Range("H16").Value = "Argentina(o)"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("N16").Value = "$/caja 50 unidades"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 290
$ws.Range("Q16").Value = 50
$ws.Range("R16").Value = "Hortaliza"
